$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 111-112, pushing the existing rows 111..168 down to 113..170
$ws.Rows("111:112").Insert()

# Row 111: new weekly record (Americana (o) / Primera)
$ws.Range("A111").Value = 2
$ws.Range("B111").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C111").Value = "Coquimbo"
$ws.Range("D111").Value = 44489
$ws.Range("E111").Value = 4
$ws.Range("F111").Value = 100112021
$ws.Range("G111").Value = "Ají"
$ws.Range("H111").Value = "Americana (o)"
$ws.Range("I111").Value = "Primera"
$ws.Range("J111").Value = 200
$ws.Range("K111").Value = 30000
$ws.Range("L111").Value = 35000
$ws.Range("M111").Value = 32500
$ws.Range("N111").Value = "`$/caja 25 kilos"
$ws.Range("O111").Value = "Provincia de Limarí"
$ws.Range("P111").Value = 1300
$ws.Range("Q111").Value = 25
$ws.Range("R111").Value = "Hortaliza"

# Row 112: new weekly record (Americana (o) / Segunda)
$ws.Range("A112").Value = 2
$ws.Range("B112").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C112").Value = "Coquimbo"
$ws.Range("D112").Value = 44489
$ws.Range("E112").Value = 4
$ws.Range("F112").Value = 100112021
$ws.Range("G112").Value = "Ají"
$ws.Range("H112").Value = "Americana (o)"
$ws.Range("I112").Value = "Segunda"
$ws.Range("J112").Value = 300
$ws.Range("K112").Value = 20000
$ws.Range("L112").Value = 25000
$ws.Range("M112").Value = 22500
$ws.Range("N112").Value = "`$/caja 25 kilos"
$ws.Range("O112").Value = "Provincia de Limarí"
$ws.Range("P112").Value = 900
$ws.Range("Q112").Value = 25
$ws.Range("R112").Value = "Hortaliza"
